# "Generate Report for Handoff" — updates the localization-status report
# after b.md moves from "Handed back" to "Ready for handoff" and a new
# handoff package is generated for both target locales.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ce95000333c970874012ebdde93376186a6d125/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/008ae04b31086f1c3d25b2772efef8dbda6316c0/e2e/b.md."

# ---- Overview sheet: row 3 is the b.md summary row ----
$ovw.Range("E3").Value = "Ready for handoff"
$ovw.Range("F3").Value = "Ready for handoff"
$ovw.Range("G3").Value = "2016-09-03 22:41:16"

# ---- zh-cn sheet: row 3 is the b.md detail row ----
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "'False"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-03 22:41:12"
$zh.Range("P3").Value = $errorDetail
$zh.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet: row 3 is the b.md detail row ----
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "'False"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-09-03 22:41:16"
$de.Range("P3").Value = $errorDetail
$de.Columns.Item(16).ColumnWidth = 39.166666666666664
